$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-ran and appended at 2026-01-16 18:38:46 JST; every row
# currently stamped with the previous run's timestamp (2026-01-16 18:29:55)
# in column A ("取得日時") picks up the new timestamp.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-01-16 18:38:46"
}
